$wb = $excel.ActiveWorkbook

# --- Sheet: final_fail ---
$ws1 = $wb.Worksheets.Item("final_fail")
$ws1.Range("A2").Value = "Total time online (min)"
$ws1.Range("B2").Value = $true
$ws1.Range("C2").Value = $true
$ws1.Range("D2").Value = $false
$ws1.Range("E2").Value = $true
$ws1.Range("F2").Value = $true
$ws1.Range("G2").Value = $true
$ws1.Range("H2").Value = $true
$ws1.Range("I2").Value = $true
$ws1.Range("J2").Value = 7

$ws1.Range("A3").Value = "Assignments viewed"
$ws1.Range("B3").Value = $true
$ws1.Range("C3").Value = $true
$ws1.Range("D3").Value = $true
$ws1.Range("E3").Value = $true
$ws1.Range("F3").Value = $true
$ws1.Range("G3").Value = $false
$ws1.Range("H3").Value = $true
$ws1.Range("I3").Value = $false
$ws1.Range("J3").Value = 6

$ws1.Range("A4").Value = "Submissions (% of course total)"
$ws1.Range("B4").Value = $true
$ws1.Range("C4").Value = $true
$ws1.Range("D4").Value = $true
$ws1.Range("E4").Value = $true
$ws1.Range("F4").Value = $true
$ws1.Range("G4").Value = $false
$ws1.Range("H4").Value = $true
$ws1.Range("I4").Value = $false
$ws1.Range("J4").Value = 6

$ws1.Range("A5").Value = "Clicks (% of course total)"
$ws1.Range("B5").Value = $true
$ws1.Range("C5").Value = $true
$ws1.Range("D5").Value = $true
$ws1.Range("E5").Value = $true
$ws1.Range("F5").Value = $true
$ws1.Range("G5").Value = $false
$ws1.Range("H5").Value = $true
$ws1.Range("I5").Value = $false
$ws1.Range("J5").Value = 6

$ws1.Range("A6").Value = "Days with no interaction"
$ws1.Range("B6").Value = $true
$ws1.Range("C6").Value = $true
$ws1.Range("D6").Value = $true
$ws1.Range("E6").Value = $true
$ws1.Range("F6").Value = $true
$ws1.Range("G6").Value = $false
$ws1.Range("H6").Value = $true
$ws1.Range("I6").Value = $false
$ws1.Range("J6").Value = 6

$ws1.Range("A7").Value = "Days with no interaction (%)"
$ws1.Range("B7").Value = $true
$ws1.Range("C7").Value = $true
$ws1.Range("D7").Value = $false
$ws1.Range("E7").Value = $true
$ws1.Range("F7").Value = $true
$ws1.Range("G7").Value = $false
$ws1.Range("H7").Value = $true
$ws1.Range("I7").Value = $false
$ws1.Range("J7").Value = 5

$ws1.Range("A8").Value = "Clicks on campus"
$ws1.Range("B8").Value = $false
$ws1.Range("C8").Value = $true
$ws1.Range("D8").Value = $true
$ws1.Range("E8").Value = $true
$ws1.Range("F8").Value = $true
$ws1.Range("G8").Value = $false
$ws1.Range("H8").Value = $true
$ws1.Range("I8").Value = $false
$ws1.Range("J8").Value = 5

$ws1.Range("A9").Value = "Largest period of inactivity (h)"
$ws1.Range("B9").Value = $true
$ws1.Range("C9").Value = $true
$ws1.Range("D9").Value = $false
$ws1.Range("E9").Value = $true
$ws1.Range("F9").Value = $true
$ws1.Range("G9").Value = $false
$ws1.Range("H9").Value = $true
$ws1.Range("I9").Value = $false
$ws1.Range("J9").Value = 5

$ws1.Range("A10").Value = "Average session duration (min)"
$ws1.Range("B10").Value = $true
$ws1.Range("C10").Value = $true
$ws1.Range("D10").Value = $true
$ws1.Range("E10").Value = $false
$ws1.Range("F10").Value = $true
$ws1.Range("G10").Value = $false
$ws1.Range("H10").Value = $true
$ws1.Range("I10").Value = $false
$ws1.Range("J10").Value = 5

$ws1.Range("A11").Value = "Start of Session 1 (%)"
$ws1.Range("B11").Value = $true
$ws1.Range("C11").Value = $true
$ws1.Range("D11").Value = $false
$ws1.Range("E11").Value = $true
$ws1.Range("F11").Value = $true
$ws1.Range("G11").Value = $false
$ws1.Range("H11").Value = $true
$ws1.Range("I11").Value = $false
$ws1.Range("J11").Value = 5

$ws1.Range("A12").Value = "Resources viewed"
$ws1.Range("B12").Value = $true
$ws1.Range("C12").Value = $true
$ws1.Range("D12").Value = $false
$ws1.Range("E12").Value = $true
$ws1.Range("F12").Value = $true
$ws1.Range("G12").Value = $false
$ws1.Range("H12").Value = $true
$ws1.Range("I12").Value = $false
$ws1.Range("J12").Value = 5

$ws1.Range("A13").Value = "Number of days"
$ws1.Range("B13").Value = $true
$ws1.Range("C13").Value = $true
$ws1.Range("D13").Value = $false
$ws1.Range("E13").Value = $true
$ws1.Range("F13").Value = $true
$ws1.Range("G13").Value = $false
$ws1.Range("H13").Value = $true
$ws1.Range("I13").Value = $false
$ws1.Range("J13").Value = 5

$ws1.Range("A14").Value = "On/off campus click ratio"
$ws1.Range("B14").Value = $true
$ws1.Range("C14").Value = $true
$ws1.Range("D14").Value = $false
$ws1.Range("E14").Value = $true
$ws1.Range("F14").Value = $true
$ws1.Range("G14").Value = $false
$ws1.Range("H14").Value = $true
$ws1.Range("I14").Value = $false
$ws1.Range("J14").Value = 5

$ws1.Range("A15").Value = "Quizzes started"
$ws1.Range("B15").Value = $true
$ws1.Range("C15").Value = $true
$ws1.Range("D15").Value = $true
$ws1.Range("E15").Value = $true
$ws1.Range("F15").Value = $false
$ws1.Range("G15").Value = $false
$ws1.Range("H15").Value = $true
$ws1.Range("I15").Value = $false
$ws1.Range("J15").Value = 5

$ws1.Range("A16").Value = "Clicks per session"
$ws1.Range("B16").Value = $true
$ws1.Range("C16").Value = $true
$ws1.Range("D16").Value = $false
$ws1.Range("E16").Value = $false
$ws1.Range("F16").Value = $true
$ws1.Range("G16").Value = $false
$ws1.Range("H16").Value = $true
$ws1.Range("I16").Value = $false
$ws1.Range("J16").Value = 4

$ws1.Range("A17").Value = "Assignments submitted"
$ws1.Range("B17").Value = $true
$ws1.Range("C17").Value = $false
$ws1.Range("D17").Value = $true
$ws1.Range("E17").Value = $true
$ws1.Range("F17").Value = $false
$ws1.Range("G17").Value = $false
$ws1.Range("H17").Value = $true
$ws1.Range("I17").Value = $false
$ws1.Range("J17").Value = 4

$ws1.Range("A18").Value = "Clicks per day"
$ws1.Range("B18").Value = $false
$ws1.Range("C18").Value = $true
$ws1.Range("D18").Value = $false
$ws1.Range("E18").Value = $true
$ws1.Range("F18").Value = $false
$ws1.Range("G18").Value = $false
$ws1.Range("H18").Value = $true
$ws1.Range("I18").Value = $false
$ws1.Range("J18").Value = 3

$ws1.Range("A19").Value = "Clicks on course"
$ws1.Range("B19").Value = $true
$ws1.Range("C19").Value = $true
$ws1.Range("D19").Value = $false
$ws1.Range("E19").Value = $false
$ws1.Range("F19").Value = $false
$ws1.Range("G19").Value = $false
$ws1.Range("H19").Value = $true
$ws1.Range("I19").Value = $false
$ws1.Range("J19").Value = 3

$ws1.Range("A20").Value = "Number of clicks"
$ws1.Range("B20").Value = $false
$ws1.Range("C20").Value = $false
$ws1.Range("D20").Value = $true
$ws1.Range("E20").Value = $true
$ws1.Range("F20").Value = $false
$ws1.Range("G20").Value = $false
$ws1.Range("H20").Value = $true
$ws1.Range("I20").Value = $false
$ws1.Range("J20").Value = 3

$ws1.Range("A21").Value = "Links viewed"
$ws1.Range("B21").Value = $false
$ws1.Range("C21").Value = $true
$ws1.Range("D21").Value = $true
$ws1.Range("E21").Value = $false
$ws1.Range("F21").Value = $false
$ws1.Range("G21").Value = $false
$ws1.Range("H21").Value = $true
$ws1.Range("I21").Value = $false
$ws1.Range("J21").Value = 3

$ws1.Range("A22").Value = "Number of sessions"
$ws1.Range("B22").Value = $false
$ws1.Range("C22").Value = $true
$ws1.Range("D22").Value = $false
$ws1.Range("E22").Value = $false
$ws1.Range("F22").Value = $true
$ws1.Range("G22").Value = $false
$ws1.Range("H22").Value = $true
$ws1.Range("I22").Value = $false
$ws1.Range("J22").Value = 3

$ws1.Range("A23").Value = "Start of Session 6 (%)"
$ws1.Range("B23").Value = $true
$ws1.Range("C23").Value = $true
$ws1.Range("D23").Value = $false
$ws1.Range("E23").Value = $false
$ws1.Range("F23").Value = $false
$ws1.Range("G23").Value = $false
$ws1.Range("H23").Value = $true
$ws1.Range("I23").Value = $false
$ws1.Range("J23").Value = 3

$ws1.Range("A24").Value = "Start of Session 3 (%)"
$ws1.Range("B24").Value = $false
$ws1.Range("C24").Value = $true
$ws1.Range("D24").Value = $false
$ws1.Range("E24").Value = $false
$ws1.Range("F24").Value = $true
$ws1.Range("G24").Value = $false
$ws1.Range("H24").Value = $true
$ws1.Range("I24").Value = $false
$ws1.Range("J24").Value = 3

$ws1.Range("A25").Value = "Start of Session 10 (%)"
$ws1.Range("B25").Value = $false
$ws1.Range("C25").Value = $true
$ws1.Range("D25").Value = $false
$ws1.Range("E25").Value = $false
$ws1.Range("F25").Value = $false
$ws1.Range("G25").Value = $false
$ws1.Range("H25").Value = $true
$ws1.Range("I25").Value = $false
$ws1.Range("J25").Value = 2

$ws1.Range("A26").Value = "Discussions viewed"
$ws1.Range("B26").Value = $false
$ws1.Range("C26").Value = $false
$ws1.Range("D26").Value = $true
$ws1.Range("E26").Value = $false
$ws1.Range("F26").Value = $false
$ws1.Range("G26").Value = $false
$ws1.Range("H26").Value = $true
$ws1.Range("I26").Value = $false
$ws1.Range("J26").Value = 2

$ws1.Range("A27").Value = "Start of Session 5 (%)"
$ws1.Range("B27").Value = $false
$ws1.Range("C27").Value = $true
$ws1.Range("D27").Value = $false
$ws1.Range("E27").Value = $false
$ws1.Range("F27").Value = $false
$ws1.Range("G27").Value = $false
$ws1.Range("H27").Value = $true
$ws1.Range("I27").Value = $false
$ws1.Range("J27").Value = 2

$ws1.Range("A28").Value = "Start of Session 2 (%)"
$ws1.Range("B28").Value = $false
$ws1.Range("C28").Value = $false
$ws1.Range("D28").Value = $false
$ws1.Range("E28").Value = $false
$ws1.Range("F28").Value = $true
$ws1.Range("G28").Value = $false
$ws1.Range("H28").Value = $true
$ws1.Range("I28").Value = $false
$ws1.Range("J28").Value = 2

$ws1.Range("A29").Value = "Clicks on forum"
$ws1.Range("B29").Value = $false
$ws1.Range("C29").Value = $true
$ws1.Range("D29").Value = $false
$ws1.Range("E29").Value = $false
$ws1.Range("F29").Value = $false
$ws1.Range("G29").Value = $false
$ws1.Range("H29").Value = $true
$ws1.Range("I29").Value = $false
$ws1.Range("J29").Value = 2

$ws1.Range("A30").Value = "Clicks on folder"
$ws1.Range("B30").Value = $false
$ws1.Range("C30").Value = $false
$ws1.Range("D30").Value = $false
$ws1.Range("E30").Value = $false
$ws1.Range("F30").Value = $false
$ws1.Range("G30").Value = $false
$ws1.Range("H30").Value = $true
$ws1.Range("I30").Value = $false
$ws1.Range("J30").Value = 1

$ws1.Range("A31").Value = "Start of Session 9 (%)"
$ws1.Range("B31").Value = $false
$ws1.Range("C31").Value = $false
$ws1.Range("D31").Value = $false
$ws1.Range("E31").Value = $false
$ws1.Range("F31").Value = $false
$ws1.Range("G31").Value = $false
$ws1.Range("H31").Value = $true
$ws1.Range("I31").Value = $false
$ws1.Range("J31").Value = 1

$ws1.Range("A32").Value = "Start of Session 8 (%)"
$ws1.Range("B32").Value = $false
$ws1.Range("C32").Value = $false
$ws1.Range("D32").Value = $false
$ws1.Range("E32").Value = $false
$ws1.Range("F32").Value = $false
$ws1.Range("G32").Value = $false
$ws1.Range("H32").Value = $true
$ws1.Range("I32").Value = $false
$ws1.Range("J32").Value = 1

$ws1.Range("A33").Value = "Forum posts"
$ws1.Range("B33").Value = $false
$ws1.Range("C33").Value = $false
$ws1.Range("D33").Value = $false
$ws1.Range("E33").Value = $false
$ws1.Range("F33").Value = $false
$ws1.Range("G33").Value = $false
$ws1.Range("H33").Value = $true
$ws1.Range("I33").Value = $false
$ws1.Range("J33").Value = 1

$ws1.Range("A34").Value = "Files downloaded"
$ws1.Range("B34").Value = $false
$ws1.Range("C34").Value = $false
$ws1.Range("D34").Value = $false
$ws1.Range("E34").Value = $false
$ws1.Range("F34").Value = $false
$ws1.Range("G34").Value = $false
$ws1.Range("H34").Value = $true
$ws1.Range("I34").Value = $false
$ws1.Range("J34").Value = 1

$ws1.Range("A35").Value = "Start of Session 7 (%)"
$ws1.Range("B35").Value = $false
$ws1.Range("C35").Value = $false
$ws1.Range("D35").Value = $false
$ws1.Range("E35").Value = $false
$ws1.Range("F35").Value = $false
$ws1.Range("G35").Value = $false
$ws1.Range("H35").Value = $true
$ws1.Range("I35").Value = $false
$ws1.Range("J35").Value = 1

$ws1.Range("A36").Value = "Start of Session 4 (%)"
$ws1.Range("B36").Value = $false
$ws1.Range("C36").Value = $false
$ws1.Range("D36").Value = $false
$ws1.Range("E36").Value = $false
$ws1.Range("F36").Value = $false
$ws1.Range("G36").Value = $false
$ws1.Range("H36").Value = $true
$ws1.Range("I36").Value = $false
$ws1.Range("J36").Value = 1

# --- Sheet: final_gifted ---
$ws2 = $wb.Worksheets.Item("final_gifted")
$ws2.Range("A2").Value = "Number of clicks"
$ws2.Range("B2").Value = $true
$ws2.Range("C2").Value = $false
$ws2.Range("D2").Value = $true
$ws2.Range("E2").Value = $true
$ws2.Range("F2").Value = $true
$ws2.Range("G2").Value = $true
$ws2.Range("H2").Value = $true
$ws2.Range("I2").Value = $true
$ws2.Range("J2").Value = 7

$ws2.Range("A3").Value = "Days with no interaction"
$ws2.Range("B3").Value = $true
$ws2.Range("C3").Value = $true
$ws2.Range("D3").Value = $false
$ws2.Range("E3").Value = $true
$ws2.Range("F3").Value = $true
$ws2.Range("G3").Value = $true
$ws2.Range("H3").Value = $true
$ws2.Range("I3").Value = $true
$ws2.Range("J3").Value = 7

$ws2.Range("A4").Value = "Clicks on campus"
$ws2.Range("B4").Value = $false
$ws2.Range("C4").Value = $false
$ws2.Range("D4").Value = $true
$ws2.Range("E4").Value = $true
$ws2.Range("F4").Value = $true
$ws2.Range("G4").Value = $true
$ws2.Range("H4").Value = $true
$ws2.Range("I4").Value = $true
$ws2.Range("J4").Value = 6

$ws2.Range("A5").Value = "Largest period of inactivity (h)"
$ws2.Range("B5").Value = $true
$ws2.Range("C5").Value = $false
$ws2.Range("D5").Value = $false
$ws2.Range("E5").Value = $true
$ws2.Range("F5").Value = $true
$ws2.Range("G5").Value = $true
$ws2.Range("H5").Value = $true
$ws2.Range("I5").Value = $true
$ws2.Range("J5").Value = 6

$ws2.Range("A6").Value = "Average session duration (min)"
$ws2.Range("B6").Value = $true
$ws2.Range("C6").Value = $true
$ws2.Range("D6").Value = $true
$ws2.Range("E6").Value = $true
$ws2.Range("F6").Value = $true
$ws2.Range("G6").Value = $false
$ws2.Range("H6").Value = $true
$ws2.Range("I6").Value = $false
$ws2.Range("J6").Value = 6

$ws2.Range("A7").Value = "Clicks (% of course total)"
$ws2.Range("B7").Value = $true
$ws2.Range("C7").Value = $true
$ws2.Range("D7").Value = $true
$ws2.Range("E7").Value = $true
$ws2.Range("F7").Value = $true
$ws2.Range("G7").Value = $false
$ws2.Range("H7").Value = $true
$ws2.Range("I7").Value = $false
$ws2.Range("J7").Value = 6

$ws2.Range("A8").Value = "Assignments viewed"
$ws2.Range("B8").Value = $true
$ws2.Range("C8").Value = $true
$ws2.Range("D8").Value = $false
$ws2.Range("E8").Value = $true
$ws2.Range("F8").Value = $false
$ws2.Range("G8").Value = $true
$ws2.Range("H8").Value = $true
$ws2.Range("I8").Value = $true
$ws2.Range("J8").Value = 6

$ws2.Range("A9").Value = "Total time online (min)"
$ws2.Range("B9").Value = $true
$ws2.Range("C9").Value = $true
$ws2.Range("D9").Value = $false
$ws2.Range("E9").Value = $true
$ws2.Range("F9").Value = $true
$ws2.Range("G9").Value = $false
$ws2.Range("H9").Value = $true
$ws2.Range("I9").Value = $false
$ws2.Range("J9").Value = 5

$ws2.Range("A10").Value = "Start of Session 1 (%)"
$ws2.Range("B10").Value = $true
$ws2.Range("C10").Value = $true
$ws2.Range("D10").Value = $false
$ws2.Range("E10").Value = $true
$ws2.Range("F10").Value = $true
$ws2.Range("G10").Value = $false
$ws2.Range("H10").Value = $true
$ws2.Range("I10").Value = $false
$ws2.Range("J10").Value = 5

$ws2.Range("A11").Value = "Clicks per session"
$ws2.Range("B11").Value = $true
$ws2.Range("C11").Value = $true
$ws2.Range("D11").Value = $false
$ws2.Range("E11").Value = $true
$ws2.Range("F11").Value = $true
$ws2.Range("G11").Value = $false
$ws2.Range("H11").Value = $true
$ws2.Range("I11").Value = $false
$ws2.Range("J11").Value = 5

$ws2.Range("A12").Value = "On/off campus click ratio"
$ws2.Range("B12").Value = $true
$ws2.Range("C12").Value = $true
$ws2.Range("D12").Value = $false
$ws2.Range("E12").Value = $true
$ws2.Range("F12").Value = $true
$ws2.Range("G12").Value = $false
$ws2.Range("H12").Value = $true
$ws2.Range("I12").Value = $false
$ws2.Range("J12").Value = 5

$ws2.Range("A13").Value = "Resources viewed"
$ws2.Range("B13").Value = $true
$ws2.Range("C13").Value = $true
$ws2.Range("D13").Value = $false
$ws2.Range("E13").Value = $true
$ws2.Range("F13").Value = $true
$ws2.Range("G13").Value = $false
$ws2.Range("H13").Value = $true
$ws2.Range("I13").Value = $false
$ws2.Range("J13").Value = 5

$ws2.Range("A14").Value = "Submissions (% of course total)"
$ws2.Range("B14").Value = $true
$ws2.Range("C14").Value = $false
$ws2.Range("D14").Value = $true
$ws2.Range("E14").Value = $false
$ws2.Range("F14").Value = $true
$ws2.Range("G14").Value = $false
$ws2.Range("H14").Value = $true
$ws2.Range("I14").Value = $false
$ws2.Range("J14").Value = 4

$ws2.Range("A15").Value = "Number of days"
$ws2.Range("B15").Value = $true
$ws2.Range("C15").Value = $false
$ws2.Range("D15").Value = $false
$ws2.Range("E15").Value = $true
$ws2.Range("F15").Value = $true
$ws2.Range("G15").Value = $false
$ws2.Range("H15").Value = $true
$ws2.Range("I15").Value = $false
$ws2.Range("J15").Value = 4

$ws2.Range("A16").Value = "Quizzes started"
$ws2.Range("B16").Value = $true
$ws2.Range("C16").Value = $false
$ws2.Range("D16").Value = $true
$ws2.Range("E16").Value = $true
$ws2.Range("F16").Value = $false
$ws2.Range("G16").Value = $false
$ws2.Range("H16").Value = $true
$ws2.Range("I16").Value = $false
$ws2.Range("J16").Value = 4

$ws2.Range("A17").Value = "Clicks on course"
$ws2.Range("B17").Value = $false
$ws2.Range("C17").Value = $true
$ws2.Range("D17").Value = $false
$ws2.Range("E17").Value = $true
$ws2.Range("F17").Value = $true
$ws2.Range("G17").Value = $false
$ws2.Range("H17").Value = $true
$ws2.Range("I17").Value = $false
$ws2.Range("J17").Value = 4

$ws2.Range("A18").Value = "Days with no interaction (%)"
$ws2.Range("B18").Value = $false
$ws2.Range("C18").Value = $false
$ws2.Range("D18").Value = $true
$ws2.Range("E18").Value = $true
$ws2.Range("F18").Value = $true
$ws2.Range("G18").Value = $false
$ws2.Range("H18").Value = $true
$ws2.Range("I18").Value = $false
$ws2.Range("J18").Value = 4

$ws2.Range("A19").Value = "Start of Session 2 (%)"
$ws2.Range("B19").Value = $false
$ws2.Range("C19").Value = $false
$ws2.Range("D19").Value = $false
$ws2.Range("E19").Value = $true
$ws2.Range("F19").Value = $true
$ws2.Range("G19").Value = $false
$ws2.Range("H19").Value = $true
$ws2.Range("I19").Value = $false
$ws2.Range("J19").Value = 3

$ws2.Range("A20").Value = "Start of Session 3 (%)"
$ws2.Range("B20").Value = $true
$ws2.Range("C20").Value = $true
$ws2.Range("D20").Value = $false
$ws2.Range("E20").Value = $false
$ws2.Range("F20").Value = $false
$ws2.Range("G20").Value = $false
$ws2.Range("H20").Value = $true
$ws2.Range("I20").Value = $false
$ws2.Range("J20").Value = 3

$ws2.Range("A21").Value = "Start of Session 7 (%)"
$ws2.Range("B21").Value = $true
$ws2.Range("C21").Value = $false
$ws2.Range("D21").Value = $false
$ws2.Range("E21").Value = $false
$ws2.Range("F21").Value = $true
$ws2.Range("G21").Value = $false
$ws2.Range("H21").Value = $true
$ws2.Range("I21").Value = $false
$ws2.Range("J21").Value = 3

$ws2.Range("A22").Value = "Clicks per day"
$ws2.Range("B22").Value = $false
$ws2.Range("C22").Value = $false
$ws2.Range("D22").Value = $true
$ws2.Range("E22").Value = $true
$ws2.Range("F22").Value = $false
$ws2.Range("G22").Value = $false
$ws2.Range("H22").Value = $true
$ws2.Range("I22").Value = $false
$ws2.Range("J22").Value = 3

$ws2.Range("A23").Value = "Start of Session 5 (%)"
$ws2.Range("B23").Value = $false
$ws2.Range("C23").Value = $false
$ws2.Range("D23").Value = $false
$ws2.Range("E23").Value = $false
$ws2.Range("F23").Value = $true
$ws2.Range("G23").Value = $false
$ws2.Range("H23").Value = $true
$ws2.Range("I23").Value = $false
$ws2.Range("J23").Value = 2

$ws2.Range("A24").Value = "Clicks on folder"
$ws2.Range("B24").Value = $false
$ws2.Range("C24").Value = $false
$ws2.Range("D24").Value = $true
$ws2.Range("E24").Value = $false
$ws2.Range("F24").Value = $false
$ws2.Range("G24").Value = $false
$ws2.Range("H24").Value = $true
$ws2.Range("I24").Value = $false
$ws2.Range("J24").Value = 2

$ws2.Range("A25").Value = "Start of Session 10 (%)"
$ws2.Range("B25").Value = $true
$ws2.Range("C25").Value = $false
$ws2.Range("D25").Value = $false
$ws2.Range("E25").Value = $false
$ws2.Range("F25").Value = $false
$ws2.Range("G25").Value = $false
$ws2.Range("H25").Value = $true
$ws2.Range("I25").Value = $false
$ws2.Range("J25").Value = 2

$ws2.Range("A26").Value = "Assignments submitted"
$ws2.Range("B26").Value = $false
$ws2.Range("C26").Value = $false
$ws2.Range("D26").Value = $true
$ws2.Range("E26").Value = $false
$ws2.Range("F26").Value = $false
$ws2.Range("G26").Value = $false
$ws2.Range("H26").Value = $true
$ws2.Range("I26").Value = $false
$ws2.Range("J26").Value = 2

$ws2.Range("A27").Value = "Start of Session 4 (%)"
$ws2.Range("B27").Value = $false
$ws2.Range("C27").Value = $false
$ws2.Range("D27").Value = $false
$ws2.Range("E27").Value = $false
$ws2.Range("F27").Value = $true
$ws2.Range("G27").Value = $false
$ws2.Range("H27").Value = $true
$ws2.Range("I27").Value = $false
$ws2.Range("J27").Value = 2

$ws2.Range("A28").Value = "Links viewed"
$ws2.Range("B28").Value = $false
$ws2.Range("C28").Value = $false
$ws2.Range("D28").Value = $true
$ws2.Range("E28").Value = $false
$ws2.Range("F28").Value = $false
$ws2.Range("G28").Value = $false
$ws2.Range("H28").Value = $true
$ws2.Range("I28").Value = $false
$ws2.Range("J28").Value = 2

$ws2.Range("A29").Value = "Start of Session 6 (%)"
$ws2.Range("B29").Value = $false
$ws2.Range("C29").Value = $false
$ws2.Range("D29").Value = $false
$ws2.Range("E29").Value = $false
$ws2.Range("F29").Value = $true
$ws2.Range("G29").Value = $false
$ws2.Range("H29").Value = $true
$ws2.Range("I29").Value = $false
$ws2.Range("J29").Value = 2

$ws2.Range("A30").Value = "Number of sessions"
$ws2.Range("B30").Value = $false
$ws2.Range("C30").Value = $false
$ws2.Range("D30").Value = $false
$ws2.Range("E30").Value = $false
$ws2.Range("F30").Value = $false
$ws2.Range("G30").Value = $false
$ws2.Range("H30").Value = $true
$ws2.Range("I30").Value = $false
$ws2.Range("J30").Value = 1

$ws2.Range("A31").Value = "Forum posts"
$ws2.Range("B31").Value = $false
$ws2.Range("C31").Value = $false
$ws2.Range("D31").Value = $false
$ws2.Range("E31").Value = $false
$ws2.Range("F31").Value = $false
$ws2.Range("G31").Value = $false
$ws2.Range("H31").Value = $true
$ws2.Range("I31").Value = $false
$ws2.Range("J31").Value = 1

$ws2.Range("A32").Value = "Files downloaded"
$ws2.Range("B32").Value = $false
$ws2.Range("C32").Value = $false
$ws2.Range("D32").Value = $false
$ws2.Range("E32").Value = $false
$ws2.Range("F32").Value = $false
$ws2.Range("G32").Value = $false
$ws2.Range("H32").Value = $true
$ws2.Range("I32").Value = $false
$ws2.Range("J32").Value = 1

$ws2.Range("A33").Value = "Discussions viewed"
$ws2.Range("B33").Value = $false
$ws2.Range("C33").Value = $false
$ws2.Range("D33").Value = $false
$ws2.Range("E33").Value = $false
$ws2.Range("F33").Value = $false
$ws2.Range("G33").Value = $false
$ws2.Range("H33").Value = $true
$ws2.Range("I33").Value = $false
$ws2.Range("J33").Value = 1

$ws2.Range("A34").Value = "Start of Session 8 (%)"
$ws2.Range("B34").Value = $false
$ws2.Range("C34").Value = $false
$ws2.Range("D34").Value = $false
$ws2.Range("E34").Value = $false
$ws2.Range("F34").Value = $false
$ws2.Range("G34").Value = $false
$ws2.Range("H34").Value = $true
$ws2.Range("I34").Value = $false
$ws2.Range("J34").Value = 1

$ws2.Range("A35").Value = "Start of Session 9 (%)"
$ws2.Range("B35").Value = $false
$ws2.Range("C35").Value = $false
$ws2.Range("D35").Value = $false
$ws2.Range("E35").Value = $false
$ws2.Range("F35").Value = $false
$ws2.Range("G35").Value = $false
$ws2.Range("H35").Value = $true
$ws2.Range("I35").Value = $false
$ws2.Range("J35").Value = 1

$ws2.Range("A36").Value = "Clicks on forum"
$ws2.Range("B36").Value = $false
$ws2.Range("C36").Value = $false
$ws2.Range("D36").Value = $false
$ws2.Range("E36").Value = $false
$ws2.Range("F36").Value = $false
$ws2.Range("G36").Value = $false
$ws2.Range("H36").Value = $true
$ws2.Range("I36").Value = $false
$ws2.Range("J36").Value = 1

